$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.414.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.933.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.55%  "
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.216.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.803"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.937.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.351.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "226.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.77%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("E31").Value = "  -5.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.21%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("E43").Value = "  -5.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.332.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("E46").Value = "  -6.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.40%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.108.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.04%  "
